$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# The "excel-range-areas" snippet (Worksheet.getRange / colorAllFormulaCells)
# was unmapped from the table, so remove its row entirely. This shifts all
# subsequent rows up by one and leaves the two now-unused shared strings to
# be dropped automatically when the workbook is saved.
$ws.Rows.Item(109).EntireRow.Delete()

# Update the active selection to reflect where the editor was last working
# (a cell further down the sheet, around the pivoted/fixed rows).
$ws.Range("O91").Select()
